$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44533
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 3600
$ws.Range("P2").Value = 3550
$ws.Range("Q2").Value = "$/kilo"
$ws.Range("R2").Value = "Región del Maule"
$ws.Range("S2").Value = 3550
$ws.Range("T2").Value = 1

# Row 3
$ws.Range("D3").Value = 44162
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("Q3").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 4667
$ws.Range("T3").Value = 1.5

# Row 4
$ws.Range("D4").Value = 44162
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 6500
$ws.Range("O4").Value = 6500
$ws.Range("P4").Value = 6500
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 4333

# Row 6
$ws.Range("D6").Value = 44159
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 6500
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6750
$ws.Range("S6").Value = 4500

# Row 8
$ws.Range("D8").Value = 44519
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 3700
$ws.Range("O8").Value = 3800
$ws.Range("P8").Value = 3750
$ws.Range("S8").Value = 3750

# Row 9
$ws.Range("D9").Value = 44537
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 5000
$ws.Range("O9").Value = 5500
$ws.Range("P9").Value = 5250
$ws.Range("R9").Value = "Región del Maule"
$ws.Range("S9").Value = 3500

# Row 10
$ws.Range("D10").Value = 44516
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 3700
$ws.Range("O10").Value = 3800
$ws.Range("P10").Value = 3750
$ws.Range("Q10").Value = "$/kilo"
$ws.Range("R10").Value = "Región del Maule"
$ws.Range("S10").Value = 3750
$ws.Range("T10").Value = 1

# Row 11
$ws.Range("D11").Value = 44176
$ws.Range("M11").Value = 300
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 5500
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value = 3667

# Row 12
$ws.Range("D12").Value = 44530
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 3600
$ws.Range("O12").Value = 3700
$ws.Range("P12").Value = 3650
$ws.Range("Q12").Value = "$/kilo"
$ws.Range("S12").Value = 3650
$ws.Range("T12").Value = 1

# Row 13
$ws.Range("D13").Value = 44523
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 3700
$ws.Range("O13").Value = 3800
$ws.Range("P13").Value = 3750
$ws.Range("Q13").Value = "$/kilo"
$ws.Range("R13").Value = "Región del Maule"
$ws.Range("S13").Value = 3750
$ws.Range("T13").Value = 1

# Row 14
$ws.Range("D14").Value = 44551
$ws.Range("M14").Value = 400
$ws.Range("N14").Value = 5000
$ws.Range("O14").Value = 5500
$ws.Range("P14").Value = 5250
$ws.Range("Q14").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("S14").Value = 3500
$ws.Range("T14").Value = 1.5

# Row 15
$ws.Range("D15").Value = 44166
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6500
$ws.Range("P15").Value = 6250
$ws.Range("Q15").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 4167
$ws.Range("T15").Value = 1.5

# Row 16
$ws.Range("D16").Value = 44544
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 5000
$ws.Range("O16").Value = 5500
$ws.Range("P16").Value = 5250
$ws.Range("Q16").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("S16").Value = 3500
$ws.Range("T16").Value = 1.5

# Row 17
$ws.Range("D17").Value = 44169
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = 5500
$ws.Range("O17").Value = 6000
$ws.Range("P17").Value = 5750
$ws.Range("S17").Value = 3833
